# Quarterly indexing bug-fix: shift each date in column A from the 1st of the
# quarter-start month to the 15th of the following month (mid-quarter date),
# leaving the revision values in column B untouched.

function Days-FromCivil($y, $m, $d) {
    if ($m -le 2) { $y2 = $y - 1 } else { $y2 = $y }
    if ($y2 -ge 0) { $era = [math]::Floor($y2 / 400) } else { $era = [math]::Floor(($y2 - 399) / 400) }
    $yoe = $y2 - $era * 400
    if ($m -gt 2) { $mAdj = $m - 3 } else { $mAdj = $m + 9 }
    $doy = [math]::Floor((153 * $mAdj + 2) / 5) + $d - 1
    $doe = $yoe * 365 + [math]::Floor($yoe / 4) - [math]::Floor($yoe / 100) + $doy
    return $era * 146097 + $doe - 719468
}

function Civil-FromDays($z) {
    $z = $z + 719468
    if ($z -ge 0) { $era = [math]::Floor($z / 146097) } else { $era = [math]::Floor(($z - 146096) / 146097) }
    $doe = $z - $era * 146097
    $yoe = [math]::Floor(($doe - [math]::Floor($doe / 1460) + [math]::Floor($doe / 36524) - [math]::Floor($doe / 146096)) / 365)
    $y = $yoe + $era * 400
    $doy = $doe - (365 * $yoe + [math]::Floor($yoe / 4) - [math]::Floor($yoe / 100))
    $mp = [math]::Floor((5 * $doy + 2) / 153)
    $d = $doy - [math]::Floor((153 * $mp + 2) / 5) + 1
    if ($mp -lt 10) { $m = $mp + 3 } else { $m = $mp - 9 }
    if ($m -le 2) { $yy = $y + 1 } else { $yy = $y }
    return @($yy, $m, $d)
}

# Excel serial-date epoch: serial 0 == 1899-12-30
$epochDays = Days-FromCivil 1899 12 30

function Serial-ToCivil($serial) {
    return Civil-FromDays ($epochDays + $serial)
}

function Civil-ToSerial($y, $m, $d) {
    return (Days-FromCivil $y $m $d) - $epochDays
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 150; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $serial = [int]$cell.Value2

    $civil = Serial-ToCivil $serial
    $y = $civil[0]
    $m = $civil[1]

    $newMonth = $m + 1
    $newYear = $y
    if ($newMonth -gt 12) {
        $newMonth = $newMonth - 12
        $newYear = $newYear + 1
    }

    $newSerial = Civil-ToSerial $newYear $newMonth 15
    $cell.Value = $newSerial
}
